$wb = $excel.ActiveWorkbook

$aiUrls = @(
    "https://monoist.atmarkit.co.jp/mn/articles/1912/06/news024.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/06/news061_2.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/06/news034.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/06/news064.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/06/news043.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/06/news048.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/05/news029.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/06/news032.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/05/news015.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/05/news058.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/05/news016.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/06/news052.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/05/news033.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/06/news070.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/05/news021.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/04/news012_2.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/05/news056.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/04/news012.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/05/news055.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1912/04/news020.html"
)

$remoteMeetingUrls = @(
    "https://monoist.atmarkit.co.jp/mn/articles/1911/27/news054.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1910/31/news018.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1910/28/news036.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1910/23/news048.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1909/26/news025.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1909/17/news024.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1909/17/news028_4.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1909/12/news059.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1908/19/news021_6.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1908/16/news029.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1908/20/news006.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1907/30/news048_3.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1907/17/news026.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1907/12/news025_3.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1906/28/news039.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1906/25/news057_3.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1906/07/news023_3.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1905/15/news067.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1905/08/news064.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1904/23/news014.html"
)

$sharingEconomyUrls = @(
    "https://monoist.atmarkit.co.jp/mn/articles/1908/21/news007.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1907/04/news018.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1906/07/news024.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1904/15/news022.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1903/29/news023.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1903/14/news005.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1902/26/news048.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1902/25/news037.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1902/06/news006.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1812/14/news059.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1812/11/news052_2.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1812/05/news020.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1812/05/news020_3.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1809/14/news019.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1806/20/news015.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1805/29/news010_3.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1710/19/news044.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1709/04/news012_3.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1708/17/news054.html",
    "https://monoist.atmarkit.co.jp/mn/articles/1707/19/news013_4.html"
)

$wsAI = $wb.Worksheets.Item("AI")
$row = 3
foreach ($u in $aiUrls) {
    $wsAI.Cells.Item($row, 1).Value = $u
    $row++
}

$wsRemoteMeeting = $wb.Worksheets.Item("遠隔会議")
$row = 3
foreach ($u in $remoteMeetingUrls) {
    $wsRemoteMeeting.Cells.Item($row, 1).Value = $u
    $row++
}

$wsSharingEconomy = $wb.Worksheets.Item("シェアリングエコノミー")
$row = 3
foreach ($u in $sharingEconomyUrls) {
    $wsSharingEconomy.Cells.Item($row, 1).Value = $u
    $row++
}
